$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 9) so the table goes from 8 data rows to 7.
$ws.Rows(9).Delete()

# --- Header row (row 1) ---
# Existing header cells B1:E1 keep "Algorithm" style (style index 1, bold/border/centered).
# Update text of existing headers and add three new styled headers F1:H1.
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"

# Copy formatting from an existing styled header cell onto the new header cells,
# then set their text (Copy+PasteSpecial brings the style without wiping text after).
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Non State std"

$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "One Sided mean"

$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "One Sided std"

# --- Data rows (row 2 .. row 8) ---
# Row 2 - LR
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8044364763692744
$ws.Range("D2").Value = 0.02884474942856064
$ws.Range("E2").Value = 0.6570731046050187
$ws.Range("F2").Value = 0.02321794113734018
$ws.Range("G2").Value = 0.7621339652781944
$ws.Range("H2").Value = 0.02597857523581242

# Row 3 - LDA
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.7969884900324209
$ws.Range("D3").Value = 0.02544947822617298
$ws.Range("E3").Value = 0.641120051727692
$ws.Range("F3").Value = 0.02214879882263832
$ws.Range("G3").Value = 0.7606103569174147
$ws.Range("H3").Value = 0.01390606589081266

# Row 4 - KNN
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.7061985748121039
$ws.Range("D4").Value = 0.03043559148890335
$ws.Range("E4").Value = 0.6011759375754986
$ws.Range("F4").Value = 0.02100305676576139
$ws.Range("G4").Value = 0.7160673259084309
$ws.Range("H4").Value = 0.01338062176378381

# Row 5 - DTREE (was CART)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7164975682661074
$ws.Range("D5").Value = 0.02956912790389477
$ws.Range("E5").Value = 0.6007505677164396
$ws.Range("F5").Value = 0.041301098137069
$ws.Range("G5").Value = 0.6692596264136759
$ws.Range("H5").Value = 0.02632842197611493

# Row 6 - RTREE
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.7164730082267388
$ws.Range("D6").Value = 0.02219425402118157
$ws.Range("E6").Value = 0.6039931730841751
$ws.Range("F6").Value = 0.02940696939236704
$ws.Range("G6").Value = 0.7306229026331386
$ws.Range("H6").Value = 0.0198524944118832

# Row 7 - XTREE
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8012535343613688
$ws.Range("D7").Value = 0.02954457540184793
$ws.Range("E7").Value = 0.6790708278274533
$ws.Range("F7").Value = 0.02042232892507847
$ws.Range("G7").Value = 0.7725136802625545
$ws.Range("H7").Value = 0.02149705229019658

# Row 8 - SVM (was NB, removed; SVM moved up from row 9)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8004103320237622
$ws.Range("D8").Value = 0.02665482451774143
$ws.Range("E8").Value = 0.6999599790278721
$ws.Range("F8").Value = 0.02229273530043684
$ws.Range("G8").Value = 0.7867708782270366
$ws.Range("H8").Value = 0.02023451298340065
